$d = $word.ActiveDocument

$d.Content.Find.Execute("91×35=", $true, $false, $false, $false, $false, $true, 1, $false, "46×23=", 2)
$d.Content.Find.Execute("62×74=", $true, $false, $false, $false, $false, $true, 1, $false, "93×18=", 2)
$d.Content.Find.Execute("68×42=", $true, $false, $false, $false, $false, $true, 1, $false, "81×39=", 2)
$d.Content.Find.Execute("67×47=", $true, $false, $false, $false, $false, $true, 1, $false, "91×54=", 2)
$d.Content.Find.Execute("26×18=", $true, $false, $false, $false, $false, $true, 1, $false, "21×68=", 2)
$d.Content.Find.Execute("24×48=", $true, $false, $false, $false, $false, $true, 1, $false, "45×23=", 2)
$d.Content.Find.Execute("74×63=", $true, $false, $false, $false, $false, $true, 1, $false, "38×77=", 2)
$d.Content.Find.Execute("55×59=", $true, $false, $false, $false, $false, $true, 1, $false, "14×55=", 2)
$d.Content.Find.Execute("97×90=", $true, $false, $false, $false, $false, $true, 1, $false, "39×65=", 2)
$d.Content.Find.Execute("79×80=", $true, $false, $false, $false, $false, $true, 1, $false, "33×48=", 2)
$d.Content.Find.Execute("81×67=", $true, $false, $false, $false, $false, $true, 1, $false, "75×94=", 2)
$d.Content.Find.Execute("74×43=", $true, $false, $false, $false, $false, $true, 1, $false, "22×46=", 2)
$d.Content.Find.Execute("75×57=", $true, $false, $false, $false, $false, $true, 1, $false, "97×84=", 2)
$d.Content.Find.Execute("82×30=", $true, $false, $false, $false, $false, $true, 1, $false, "16×64=", 2)
$d.Content.Find.Execute("27×55=", $true, $false, $false, $false, $false, $true, 1, $false, "51×44=", 2)
$d.Content.Find.Execute("99×54=", $true, $false, $false, $false, $false, $true, 1, $false, "34×86=", 2)
$d.Content.Find.Execute("98×35=", $true, $false, $false, $false, $false, $true, 1, $false, "94×52=", 2)
$d.Content.Find.Execute("90×79=", $true, $false, $false, $false, $false, $true, 1, $false, "13×58=", 2)
$d.Content.Find.Execute("30×27=", $true, $false, $false, $false, $false, $true, 1, $false, "72×30=", 2)
$d.Content.Find.Execute("69×38=", $true, $false, $false, $false, $false, $true, 1, $false, "81×89=", 2)
$d.Content.Find.Execute("15×39=", $true, $false, $false, $false, $false, $true, 1, $false, "75×21=", 2)
$d.Content.Find.Execute("92×92=", $true, $false, $false, $false, $false, $true, 1, $false, "78×96=", 2)
$d.Content.Find.Execute("40×13=", $true, $false, $false, $false, $false, $true, 1, $false, "27×48=", 2)
$d.Content.Find.Execute("33×32=", $true, $false, $false, $false, $false, $true, 1, $false, "14×21=", 2)
$d.Content.Find.Execute("29×51=", $true, $false, $false, $false, $false, $true, 1, $false, "42×18=", 2)
